$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data point was added to the top of this block (row 358),
# pushing the existing rows 358-466 down by one (to 359-467).
$ws.Rows.Item(358).Insert()

# Populate the newly inserted row 358 with the new record.
$ws.Range("A358").Value = 10
$ws.Range("B358").Value = "Vega Modelo de Temuco"
$ws.Range("C358").Value = "La Araucanía"
$ws.Range("D358").Value = 44985
$ws.Range("E358").Value = 9
$ws.Range("F358").Value = 100112009
$ws.Range("G358").Value = "Acelga"
$ws.Range("H358").Value = "Sin especificar"
$ws.Range("I358").Value = "Primera"
$ws.Range("J358").Value = 50
$ws.Range("K358").Value = 8000
$ws.Range("L358").Value = 8000
$ws.Range("M358").Value = 8000
$ws.Range("N358").Value = "$/docena de atados (12 kilos)"
$ws.Range("O358").Value = "Provincia de Cautín"
$ws.Range("P358").Value = 667
$ws.Range("Q358").Value = 12
$ws.Range("R358").Value = "Hortaliza"
